$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 12 values: Current_Ct_Day, Current_Pct_Ct, Current_Ct_Tokens, Current_Pct_Tokens
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0.001
$ws.Range("K12").Value = 695
$ws.Range("L12").Value = 0.003475
